$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells stay text (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '60.706.85'
$ws.Range("E2").Value = '  -1.65%  '
$ws.Range("D3").Value = '2.906.92'
$ws.Range("E3").Value = '  -2.47%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '528.08'
$ws.Range("E5").Value = '  -3.10%  '
$ws.Range("D6").Value = '143.44'
$ws.Range("E6").Value = '  -6.33%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.554'
$ws.Range("E8").Value = '  -3.54%  '
$ws.Range("D9").Value = '2.915.72'
$ws.Range("E9").Value = '  -2.51%  '
$ws.Range("E10").Value = '  -5.17%  '
$ws.Range("D11").Value = '6.02'
$ws.Range("E11").Value = '  -2.16%  '
$ws.Range("E12").Value = '  -2.66%  '
$ws.Range("D13").Value = '3.418.61'
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("D14").Value = '0.128'
$ws.Range("E14").Value = '  +2.47%  '
$ws.Range("D15").Value = '60.679.21'
$ws.Range("E15").Value = '  -1.79%  '
$ws.Range("D16").Value = '22.84'
$ws.Range("E16").Value = '  -4.02%  '
$ws.Range("D17").Value = '2.915.83'
$ws.Range("E17").Value = '  -2.25%  '
$ws.Range("D18").Value = '0.0000141'
$ws.Range("E18").Value = '  -4.32%  '
$ws.Range("E19").Value = '  -2.84%  '
$ws.Range("D20").Value = '11.71'
$ws.Range("E20").Value = '  -2.79%  '
$ws.Range("D21").Value = '361.10'
$ws.Range("E21").Value = '  -5.81%  '
$ws.Range("D22").Value = '6.65'
$ws.Range("E22").Value = '  -0.95%  '
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("D24").Value = '5.69'
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("D25").Value = '64.92'
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("D26").Value = '0.455'
$ws.Range("E26").Value = '  -3.55%  '
$ws.Range("E27").Value = '  -4.17%  '
$ws.Range("E28").Value = '  +1.53%  '
$ws.Range("D29").Value = '7.90'
$ws.Range("E29").Value = '  -4.95%  '
$ws.Range("D30").Value = '0.0₃0846'
$ws.Range("E30").Value = '  -10.75%  '
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("D32").Value = '1.69'
$ws.Range("E32").Value = '  -2.29%  '
$ws.Range("D33").Value = '19.83'
$ws.Range("E33").Value = '  -3.37%  '
$ws.Range("D34").Value = '150.08'
$ws.Range("E34").Value = '  -6.72%  '
$ws.Range("D35").Value = '4.36'
$ws.Range("E35").Value = '  -7.15%  '
$ws.Range("D36").Value = '5.58'
$ws.Range("E36").Value = '  -6.77%  '
$ws.Range("E37").Value = '  -6.94%  '
$ws.Range("D38").Value = '1.20'
$ws.Range("E38").Value = '  -5.73%  '
$ws.Range("D39").Value = '37.95'
$ws.Range("E39").Value = '  +1.44%  '
$ws.Range("E40").Value = '  -4.96%  '
$ws.Range("D41").Value = '3.72'
$ws.Range("E41").Value = '  -5.35%  '
$ws.Range("D42").Value = '2.295.99'
$ws.Range("E42").Value = '  -4.88%  '
$ws.Range("D43").Value = '0.650'
$ws.Range("E43").Value = '  -2.87%  '
$ws.Range("E44").Value = '  -1.53%  '
$ws.Range("D45").Value = '20.47'
$ws.Range("E45").Value = '  -8.16%  '
$ws.Range("D46").Value = '0.998'
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("D47").Value = '4.97'
$ws.Range("E47").Value = '  -2.78%  '
$ws.Range("E48").Value = '  -4.05%  '
$ws.Range("E49").Value = '  -1.22%  '
$ws.Range("D50").Value = '0.0922'
$ws.Range("E50").Value = '  -3.26%  '
$ws.Range("D51").Value = '249.40'
$ws.Range("E51").Value = '  -8.33%  '
